$wb = $excel.ActiveWorkbook

# --- Insert the new "Codes" worksheet right after "LoginCredentials" ---
$login = $wb.Worksheets.Item("LoginCredentials")
$codes = $wb.Worksheets.Add($null, $login)
$codes.Name = "Codes"

# Column widths (characters) - approximate the original authored widths
$codes.Columns.Item(1).ColumnWidth = 14.57
$codes.Columns.Item(2).ColumnWidth = 21.43

# Header row
$codes.Range("A1").Value = "TestCaseCodesString"
$codes.Range("B1").Value = "AlreadyUsedClassName"
$codes.Range("C1").Value = "CodesId"
$codes.Range("D1").Value = "CodesName"
$codes.Range("E1").Value = "CodesDefinition"
$codes.Range("F1").Value = "CritId"
$codes.Range("G1").Value = "OtherCodeClassFromCodeClass"
$codes.Range("H1").Value = "IndxCard"

# Data row
$codes.Range("A2").Value = "Test case code"
$codes.Range("B2").Value = "Test case sp code"
$codes.Range("C2").Value = "C"
$codes.Range("D2").Value = "Test"
$codes.Range("E2").Value = "Summary test"
$codes.Range("F2").Value = "Critical"
$codes.Range("G2").Value = "Test doc class"
$codes.Range("H2").Value = "IndxCard"

# --- LoginCredentials: drop the last data row (the "reports"/"Reports Test" row) ---
$login.Rows.Item(14).Delete()

# --- Selections / active-tab bookkeeping ---
[void]$login.Range("E19").Select()
[void]$codes.Range("G11").Select()
[void]$codes.Select()
